# The "Prix Spot" sheet gets a brand-new daily column ("16-dec") inserted
# right before the existing "01-oct." column (currently column ES). Excel's
# EntireColumn.Insert() shifts everything from ES:FW to ET:FX (and the
# header formatting/styles move along with it), growing the used range
# from A1:FW25 to A1:FX25.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Prix Spot")
$ws.Range("ES1").EntireColumn.Insert()

# Populate the freshly inserted column: header label in row 1, and the
# usual "-" placeholder (no price yet) for every data row 2..25.
$ws.Range("ES1").Value = "16-dec"
$ws.Range("ES2:ES25").Value = "-"

# The "Gaz" sheet's last two rows (gas price for 2025-12-13 / 2025-12-14)
# were updated with a refreshed price.
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("B178").Value = 25.93
$wsGaz.Range("B179").Value = 25.93
